{"js": "const pairs = [\n  [\"75\u00d756=\", \"23\u00d773=\"],\n  [\"17\u00d734=\", \"69\u00d785=\"],\n  [\"56\u00d762=\", \"32\u00d778=\"],\n  [\"66\u00d785=\", \"18\u00d742=\"],\n  [\"25\u00d715=\", \"11\u00d730=\"],\n  [\"55\u00d779=\", \"99\u00d720=\"],\n  [\"50\u00d736=\", \"42\u00d725=\"],\n  [\"54\u00d713=\", \"59\u00d726=\"],\n  [\"53\u00d748=\", \"82\u00d778=\"],\n  [\"37\u00d719=\", \"19\u00d797=\"],\n  [\"56\u00d799=\", \"18\u00d775=\"],\n  [\"75\u00d786=\", \"44\u00d739=\"],\n  [\"97\u00d782=\", \"57\u00d718=\"],\n  [\"93\u00d720=\", \"30\u00d772=\"],\n  [\"42\u00d717=\", \"52\u00d721=\"],\n  [\"43\u00d736=\", \"19\u00d767=\"],\n  [\"17\u00d790=\", \"95\u00d717=\"],\n  [\"52\u00d781=\", \"60\u00d799=\"],\n  [\"70\u00d773=\", \"99\u00d793=\"],\n  [\"72\u00d730=\", \"53\u00d746=\"],\n  [\"47\u00d717=\", \"24\u00d722=\"],\n  [\"40\u00d793=\", \"82\u00d712=\"],\n  [\"48\u00d768=\", \"14\u00d731=\"],\n  [\"46\u00d790=\", \"39\u00d751=\"],\n  [\"45\u00d761=\", \"96\u00d724=\"],\n  [\"92\u00d741=\", \"82\u00d775=\"],\n  [\"14\u00d752=\", \"27\u00d728=\"],\n  [\"32\u00d719=\", \"83\u00d794=\"],\n  [\"98\u00d778=\", \"62\u00d749=\"],\n  [\"72\u00d774=\", \"61\u00d734=\"],\n  [\"13\u00d752=\", \"70\u00d716=\"],\n  [\"44\u00d719=\", \"97\u00d728=\"],\n  [\"38\u00d746=\", \"84\u00d745=\"],\n  [\"90\u00d727=\", \"75\u00d776=\"],\n  [\"56\u00d780=\", \"72\u00d741=\"],\n  [\"60\u00d741=\", \"50\u00d792=\"],\n  [\"90\u00d770=\", \"77\u00d774=\"],\n  [\"19\u00d792=\", \"20\u00d740=\"],\n  [\"60\u00d768=\", \"73\u00d712=\"],\n  [\"54\u00d774=\", \"71\u00d789=\"],\n  [\"50\u00d739=\", \"71\u00d793=\"],\n  [\"70\u00d763=\", \"44\u00d734=\"],\n  [\"19\u00d761=\", \"22\u00d714=\"],\n  [\"22\u00d775=\", \"54\u00d717=\"],\n  [\"99\u00d737=\", \"75\u00d728=\"],\n  [\"71\u00d754=\", \"72\u00d764=\"],\n  [\"15\u00d711=\", \"25\u00d755=\"],\n  [\"25\u00d750=\", \"72\u00d767=\"],\n  [\"84\u00d775=\", \"12\u00d749=\"],\n  [\"47\u00d795=\", \"60\u00d713=\"],\n  [\"57\u00d784=\", \"26\u00d777=\"],\n  [\"16\u00d750=\", \"100\u00d782=\"],\n  [\"98\u00d782=\", \"12\u00d720=\"],\n  [\"94\u00d753=\", \"81\u00d744=\"],\n  [\"55\u00d720=\", \"87\u00d774=\"],\n  [\"70\u00d735=\", \"39\u00d792=\"],\n  [\"25\u00d771=\", \"62\u00d774=\"],\n  [\"36\u00d741=\", \"21\u00d7100=\"],\n  [\"88\u00d727=\", \"38\u00d744=\"],\n  [\"76\u00d730=\", \"47\u00d774=\"],\n  [\"67\u00d753=\", \"69\u00d755=\"],\n  [\"50\u00d743=\", \"50\u00d740=\"],\n  [\"94\u00d791=\", \"43\u00d785=\"],\n  [\"62\u00d721=\", \"92\u00d724=\"],\n  [\"83\u00d766=\", \"37\u00d784=\"],\n  [\"70\u00d750=\", \"55\u00d750=\"],\n  [\"79\u00d715=\", \"62\u00d713=\"],\n  [\"58\u00d716=\", \"45\u00d716=\"],\n  [\"48\u00d760=\", \"57\u00d713=\"],\n  [\"55\u00d772=\", \"36\u00d786=\"],\n  [\"62\u00d752=\", \"65\u00d748=\"],\n  [\"10\u00d734=\", \"70\u00d751=\"],\n  [\"69\u00d766=\", \"65\u00d743=\"],\n  [\"53\u00d779=\", \"20\u00d752=\"],\n  [\"16\u00d793=\", \"73\u00d780=\"],\n  [\"90\u00d739=\", \"19\u00d759=\"],\n  [\"78\u00d739=\", \"45\u00d773=\"],\n  [\"24\u00d782=\", \"22\u00d738=\"],\n  [\"10\u00d764=\", \"14\u00d771=\"],\n  [\"14\u00d738=\", \"63\u00d746=\"],\n  [\"45\u00d789=\", \"99\u00d782=\"],\n  [\"87\u00d740=\", \"97\u00d734=\"],\n  [\"40\u00d725=\", \"12\u00d736=\"],\n  [\"46\u00d782=\", \"53\u00d795=\"],\n  [\"63\u00d741=\", \"70\u00d794=\"],\n  [\"91\u00d755=\", \"22\u00d719=\"],\n  [\"32\u00d723=\", \"55\u00d795=\"],\n  [\"40\u00d717=\", \"97\u00d745=\"],\n  [\"55\u00d746=\", \"88\u00d789=\"],\n  [\"51\u00d789=\", \"32\u00d771=\"],\n  [\"49\u00d730=\", \"31\u00d790=\"],\n  [\"66\u00d718=\", \"39\u00d759=\"],\n  [\"41\u00d747=\", \"96\u00d761=\"],\n  [\"59\u00d764=\", \"69\u00d748=\"],\n  [\"76\u00d715=\", \"83\u00d781=\"],\n  [\"51\u00d753=\", \"11\u00d761=\"],\n  [\"47\u00d776=\", \"70\u00d731=\"],\n  [\"44\u00d748=\", \"23\u00d777=\"],\n  [\"74\u00d793=\", \"21\u00d759=\"],\n  [\"34\u00d786=\", \"33\u00d785=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$pairs = @(\n  ,@('75\u00d756=', '23\u00d773=')\n  ,@('17\u00d734=', '69\u00d785=')\n  ,@('56\u00d762=', '32\u00d778=')\n  ,@('66\u00d785=', '18\u00d742=')\n  ,@('25\u00d715=', '11\u00d730=')\n  ,@('55\u00d779=', '99\u00d720=')\n  ,@('50\u00d736=', '42\u00d725=')\n  ,@('54\u00d713=', '59\u00d726=')\n  ,@('53\u00d748=', '82\u00d778=')\n  ,@('37\u00d719=', '19\u00d797=')\n  ,@('56\u00d799=', '18\u00d775=')\n  ,@('75\u00d786=', '44\u00d739=')\n  ,@('97\u00d782=', '57\u00d718=')\n  ,@('93\u00d720=', '30\u00d772=')\n  ,@('42\u00d717=', '52\u00d721=')\n  ,@('43\u00d736=', '19\u00d767=')\n  ,@('17\u00d790=', '95\u00d717=')\n  ,@('52\u00d781=', '60\u00d799=')\n  ,@('70\u00d773=', '99\u00d793=')\n  ,@('72\u00d730=', '53\u00d746=')\n  ,@('47\u00d717=', '24\u00d722=')\n  ,@('40\u00d793=', '82\u00d712=')\n  ,@('48\u00d768=', '14\u00d731=')\n  ,@('46\u00d790=', '39\u00d751=')\n  ,@('45\u00d761=', '96\u00d724=')\n  ,@('92\u00d741=', '82\u00d775=')\n  ,@('14\u00d752=', '27\u00d728=')\n  ,@('32\u00d719=', '83\u00d794=')\n  ,@('98\u00d778=', '62\u00d749=')\n  ,@('72\u00d774=', '61\u00d734=')\n  ,@('13\u00d752=', '70\u00d716=')\n  ,@('44\u00d719=', '97\u00d728=')\n  ,@('38\u00d746=', '84\u00d745=')\n  ,@('90\u00d727=', '75\u00d776=')\n  ,@('56\u00d780=', '72\u00d741=')\n  ,@('60\u00d741=', '50\u00d792=')\n  ,@('90\u00d770=', '77\u00d774=')\n  ,@('19\u00d792=', '20\u00d740=')\n  ,@('60\u00d768=', '73\u00d712=')\n  ,@('54\u00d774=', '71\u00d789=')\n  ,@('50\u00d739=', '71\u00d793=')\n  ,@('70\u00d763=', '44\u00d734=')\n  ,@('19\u00d761=', '22\u00d714=')\n  ,@('22\u00d775=', '54\u00d717=')\n  ,@('99\u00d737=', '75\u00d728=')\n  ,@('71\u00d754=', '72\u00d764=')\n  ,@('15\u00d711=', '25\u00d755=')\n  ,@('25\u00d750=', '72\u00d767=')\n  ,@('84\u00d775=', '12\u00d749=')\n  ,@('47\u00d795=', '60\u00d713=')\n  ,@('57\u00d784=', '26\u00d777=')\n  ,@('16\u00d750=', '100\u00d782=')\n  ,@('98\u00d782=', '12\u00d720=')\n  ,@('94\u00d753=', '81\u00d744=')\n  ,@('55\u00d720=', '87\u00d774=')\n  ,@('70\u00d735=', '39\u00d792=')\n  ,@('25\u00d771=', '62\u00d774=')\n  ,@('36\u00d741=', '21\u00d7100=')\n  ,@('88\u00d727=', '38\u00d744=')\n  ,@('76\u00d730=', '47\u00d774=')\n  ,@('67\u00d753=', '69\u00d755=')\n  ,@('50\u00d743=', '50\u00d740=')\n  ,@('94\u00d791=', '43\u00d785=')\n  ,@('62\u00d721=', '92\u00d724=')\n  ,@('83\u00d766=', '37\u00d784=')\n  ,@('70\u00d750=', '55\u00d750=')\n  ,@('79\u00d715=', '62\u00d713=')\n  ,@('58\u00d716=', '45\u00d716=')\n  ,@('48\u00d760=', '57\u00d713=')\n  ,@('55\u00d772=', '36\u00d786=')\n  ,@('62\u00d752=', '65\u00d748=')\n  ,@('10\u00d734=', '70\u00d751=')\n  ,@('69\u00d766=', '65\u00d743=')\n  ,@('53\u00d779=', '20\u00d752=')\n  ,@('16\u00d793=', '73\u00d780=')\n  ,@('90\u00d739=', '19\u00d759=')\n  ,@('78\u00d739=', '45\u00d773=')\n  ,@('24\u00d782=', '22\u00d738=')\n  ,@('10\u00d764=', '14\u00d771=')\n  ,@('14\u00d738=', '63\u00d746=')\n  ,@('45\u00d789=', '99\u00d782=')\n  ,@('87\u00d740=', '97\u00d734=')\n  ,@('40\u00d725=', '12\u00d736=')\n  ,@('46\u00d782=', '53\u00d795=')\n  ,@('63\u00d741=', '70\u00d794=')\n  ,@('91\u00d755=', '22\u00d719=')\n  ,@('32\u00d723=', '55\u00d795=')\n  ,@('40\u00d717=', '97\u00d745=')\n  ,@('55\u00d746=', '88\u00d789=')\n  ,@('51\u00d789=', '32\u00d771=')\n  ,@('49\u00d730=', '31\u00d790=')\n  ,@('66\u00d718=', '39\u00d759=')\n  ,@('41\u00d747=', '96\u00d761=')\n  ,@('59\u00d764=', '69\u00d748=')\n  ,@('76\u00d715=', '83\u00d781=')\n  ,@('51\u00d753=', '11\u00d761=')\n  ,@('47\u00d776=', '70\u00d731=')\n  ,@('44\u00d748=', '23\u00d777=')\n  ,@('74\u00d793=', '21\u00d759=')\n  ,@('34\u00d786=', '33\u00d785=')\n)\n\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}"}
